# Regenerate merged AHB files
#
# The "Vorgang"/segment-header rows (one per top-level segment group, e.g.
# row 13 = "Nachrichtendatum") get re-styled to the grey "section header"
# look (fill D9D9D9, column B bold) that row 9 already uses as a template,
# and — like every other detail row under that header — the "L" column's
# "AENDERUNG" flag (style 4: bold/gold/centered, text "AENDERUNG") is
# replaced by an empty grey/centered cell (style 5), which row 3's "L"
# column already uses as a template.
#
# Rather than poke style indices directly (which is brittle across engines),
# we copy-format from existing cells that already carry the exact target
# style so the workbook's existing style catalog is reused verbatim.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose entire A:V span is restyled to the grey "section header" look
# (column B bold-grey, everything else plain-grey; L cleared/greyed).
$fullStyleRows = @(13, 17, 23, 27, 34, 40, 71, 104)

# Rows where only the "L" column's AENDERUNG flag is removed (style 4 -> 5,
# value cleared), the rest of the row is untouched.
$lOnlyRows = @(14, 15, 16, 18, 19, 20, 21, 22, 24, 25, 26, 28, 29, 30, 31, 32, 33, 35, 36, 38, 39, 41, 42, 43, 103, 108, 112, 113)

$xlPasteFormats = -4122

foreach ($r in $fullStyleRows) {
    $ws.Range("A9:V9").Copy() | Out-Null
    $ws.Range("A" + $r + ":V" + $r).PasteSpecial($xlPasteFormats) | Out-Null
    $excel.CutCopyMode = $false
    $ws.Range("L" + $r).ClearContents()
}

foreach ($r in $lOnlyRows) {
    $ws.Range("L3").Copy() | Out-Null
    $ws.Range("L" + $r).PasteSpecial($xlPasteFormats) | Out-Null
    $excel.CutCopyMode = $false
    $ws.Range("L" + $r).ClearContents()
}

Write-Host "Applied AHB regeneration restyle to $($fullStyleRows.Count) header rows and $($lOnlyRows.Count) L-only rows."
